# Updating filtered feeds from workflow
# Appends two new feed rows (37 and 38) to the "Filtered Feeds" sheet,
# mirroring the data added by the upstream feed-filtering workflow.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 37: GenomeWeb article about Caris Life Sciences -------------------
$link37 = "https://www.genomeweb.com/cancer/caris-life-sciences-eyes-wider-clinical-adoption-tests-new-validation-data"

$ws.Range("A37").Value = $link37
$ws.Range("B37").Value = "CDx"
$ws.Range("C37").Value = "FDA Approves Agilent CDx Assay for Use With Bristol Myers Squibb Immunotherapies"

$ws.Hyperlinks.Add($ws.Range("A37"), $link37)
$ws.Range("A37").Style = $ws.Range("A36").Style

# --- Row 38: 360Dx article about Caris Life Sciences ------------------------
$link38 = "https://www.360dx.com/cancer/caris-life-sciences-eyes-wider-clinical-adoption-tests-new-validation-data"

$ws.Range("A38").Value = $link38
$ws.Range("B38").Value = "CDx"
$ws.Range("C38").Value = "Caris Life Sciences Eyes Wider Clinical Adoption of Tests With New Validation Data"

$ws.Hyperlinks.Add($ws.Range("A38"), $link38)
$ws.Range("A38").Style = $ws.Range("A36").Style
